# Weekly update: insert two new daily price records (rows 166 and 167) for
# "Berenjena" at the top of the historical block, pushing the previously
# existing rows 166-200 down to rows 168-202 (dimension grows from
# A1:R200 to A1:R202).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 166; this shifts the existing rows 166-200
# down to 168-202 and extends the used range accordingly.
$ws.Range("A166:A167").EntireRow.Insert()

# ---- New row 166 ----
$ws.Range("A166").Value = 9
$ws.Range("B166").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C166").Value = "Metropolitana"
$ws.Range("D166").Value = 44522
$ws.Range("E166").Value = 13
$ws.Range("F166").Value = 100112001
$ws.Range("G166").Value = "Berenjena"
$ws.Range("H166").Value = "Sin especificar"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 160
$ws.Range("K166").Value = 9000
$ws.Range("L166").Value = 10000
$ws.Range("M166").Value = 9500
$ws.Range("N166").Value = "`$/caja 50 unidades"
$ws.Range("O166").Value = "Región de Arica y Parinacota"
$ws.Range("P166").Value = 190
$ws.Range("Q166").Value = 50
$ws.Range("R166").Value = "Hortaliza"

# ---- New row 167 ----
$ws.Range("A167").Value = 9
$ws.Range("B167").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C167").Value = "Metropolitana"
$ws.Range("D167").Value = 44522
$ws.Range("E167").Value = 13
$ws.Range("F167").Value = 100112001
$ws.Range("G167").Value = "Berenjena"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Segunda"
$ws.Range("J167").Value = 61
$ws.Range("K167").Value = 7000
$ws.Range("L167").Value = 7000
$ws.Range("M167").Value = 7000
$ws.Range("N167").Value = "`$/caja 100 unidades"
$ws.Range("O167").Value = "Región de Arica y Parinacota"
$ws.Range("P167").Value = 70
$ws.Range("Q167").Value = 100
$ws.Range("R167").Value = "Hortaliza"
